# Weekly fruit/vegetable price update: insert a new daily record as row 213,
# pushing the existing rows 213-230 down to 214-231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 213 (shifts 213:230 -> 214:231).
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new weekly record.
$ws.Cells.Item(213, 1).Value2 = 1
$ws.Cells.Item(213, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(213, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(213, 4).Value2 = 44578
$ws.Cells.Item(213, 5).Value2 = 15
$ws.Cells.Item(213, 6).Value2 = 100114013
$ws.Cells.Item(213, 7).Value2 = "Zanahoria"
$ws.Cells.Item(213, 8).Value2 = "Sin especificar"
$ws.Cells.Item(213, 9).Value2 = "Primera"
$ws.Cells.Item(213, 10).Value2 = 70
$ws.Cells.Item(213, 11).Value2 = 22000
$ws.Cells.Item(213, 12).Value2 = 23000
$ws.Cells.Item(213, 13).Value2 = 22500
$ws.Cells.Item(213, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(213, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(213, 16).Value2 = 900
$ws.Cells.Item(213, 17).Value2 = 25
$ws.Cells.Item(213, 18).Value2 = "Hortaliza"
